# Natmi LR-pairs recompute (Apoe-Ldlr): ligand/receptor expressing-cell counts
# went from 1 to 3 per Dr Hou's advice, which cascades into the dependent
# expression/specificity/edge-weight statistics in columns E,G-K,M-T for
# rows 2-17. Values below match the recomputed workbook exactly.
$data = @{
    2 = @{ "E"="3"; "G"="29.32133366666666"; "H"="87.964001"; "I"="0.006401919837078288"; "J"="0.006401919837078288"; "K"="3"; "M"="3.021894333333333"; "N"="9.065683"; "O"="0.1464771679819186"; "P"="0.1464771679819185"; "Q"="88.60597205307589"; "R"="797.453748477683"; "S"="0.000937735087382493"; "T"="0.000937735087382493" }
    3 = @{ "E"="3"; "G"="29.32133366666666"; "H"="87.964001"; "I"="0.006401919837078288"; "J"="0.006401919837078288"; "K"="3"; "M"="5.799695333333333"; "N"="17.399086"; "O"="0.2811226515149324"; "P"="0.2811226515149324"; "Q"="170.0548020336762"; "R"="1530.493218303086"; "S"="0.001799724679385492"; "T"="0.001799724679385493" }
    4 = @{ "E"="3"; "G"="29.32133366666666"; "H"="87.964001"; "I"="0.006401919837078288"; "J"="0.006401919837078288"; "K"="3"; "M"="6.259039333333333"; "N"="18.777118"; "O"="0.303387959572633"; "P"="0.303387959572633"; "Q"="183.5233807254575"; "R"="1651.710426529118"; "S"="0.001942265396718745"; "T"="0.001942265396718745" }
    5 = @{ "E"="3"; "G"="29.32133366666666"; "H"="87.964001"; "I"="0.006401919837078288"; "J"="0.006401919837078288"; "K"="3"; "M"="5.549851333333334"; "N"="16.649554"; "O"="0.2690122209305161"; "P"="0.2690122209305161"; "Q"="162.7290427450615"; "R"="1464.561384705554"; "S"="0.001722194673591558"; "T"="0.001722194673591558" }
    6 = @{ "E"="3"; "G"="45.524413"; "H"="136.573239"; "I"="0.009939644832300594"; "J"="0.009939644832300592"; "K"="3"; "M"="3.021894333333333"; "N"="9.065683"; "O"="0.1464771679819186"; "P"="0.1464771679819185"; "Q"="137.5699656730264"; "R"="1238.129691057237"; "S"="0.001455931025781503"; "T"="0.001455931025781502" }
    7 = @{ "E"="3"; "G"="45.524413"; "H"="136.573239"; "I"="0.009939644832300594"; "J"="0.009939644832300592"; "K"="3"; "M"="5.799695333333333"; "N"="17.399086"; "O"="0.2811226515149324"; "P"="0.2811226515149324"; "Q"="264.0277256288393"; "R"="2376.249530659554"; "S"="0.002794259310373039"; "T"="0.002794259310373038" }
    8 = @{ "E"="3"; "G"="45.524413"; "H"="136.573239"; "I"="0.009939644832300594"; "J"="0.009939644832300592"; "K"="3"; "M"="6.259039333333333"; "N"="18.777118"; "O"="0.303387959572633"; "P"="0.303387959572633"; "Q"="284.9390915939113"; "R"="2564.451824345202"; "S"="0.003015568564548343"; "T"="0.003015568564548342" }
    9 = @{ "E"="3"; "G"="45.524413"; "H"="136.573239"; "I"="0.009939644832300594"; "J"="0.009939644832300592"; "K"="3"; "M"="5.549851333333334"; "N"="16.649554"; "O"="0.2690122209305161"; "P"="0.2690122209305161"; "Q"="252.6537241872674"; "R"="2273.883517685406"; "S"="0.00267388593159771"; "T"="0.00267388593159771" }
    10 = @{ "E"="3"; "G"="4438.215250666667"; "H"="13314.645752"; "I"="0.9690247577915309"; "J"="0.9690247577915307"; "K"="3"; "M"="3.021894333333333"; "N"="9.065683"; "O"="0.1464771679819186"; "P"="0.1464771679819185"; "Q"="13411.81751610318"; "R"="120706.3576449286"; "S"="0.141940002225668"; "T"="0.141940002225668" }
    11 = @{ "E"="3"; "G"="4438.215250666667"; "H"="13314.645752"; "I"="0.9690247577915309"; "J"="0.9690247577915307"; "K"="3"; "M"="5.799695333333333"; "N"="17.399086"; "O"="0.2811226515149324"; "P"="0.2811226515149324"; "Q"="25740.2962776203"; "R"="231662.6664985827"; "S"="0.2724148092939703"; "T"="0.2724148092939703" }
    12 = @{ "E"="3"; "G"="4438.215250666667"; "H"="13314.645752"; "I"="0.9690247577915309"; "J"="0.9690247577915307"; "K"="3"; "M"="6.259039333333333"; "N"="18.777118"; "O"="0.303387959572633"; "P"="0.303387959572633"; "Q"="27778.96382372253"; "R"="250010.6744135027"; "S"="0.2939904440417374"; "T"="0.2939904440417374" }
    13 = @{ "E"="3"; "G"="4438.215250666667"; "H"="13314.645752"; "I"="0.9690247577915309"; "J"="0.9690247577915307"; "K"="3"; "M"="5.549851333333334"; "N"="16.649554"; "O"="0.2690122209305161"; "P"="0.2690122209305161"; "Q"="24631.43482653274"; "R"="221682.9134387946"; "S"="0.2606795022301552"; "T"="0.2606795022301551" }
    14 = @{ "E"="3"; "G"="67.02347933333333"; "H"="201.070438"; "I"="0.01463367753909034"; "J"="0.01463367753909034"; "K"="3"; "M"="3.021894333333333"; "N"="9.065683"; "O"="0.1464771679819186"; "P"="0.1464771679819185"; "Q"="202.5378723976838"; "R"="1822.840851579154"; "S"="0.002143499643086564"; "T"="0.002143499643086564" }
    15 = @{ "E"="3"; "G"="67.02347933333333"; "H"="201.070438"; "I"="0.01463367753909034"; "J"="0.01463367753909034"; "K"="3"; "M"="5.799695333333333"; "N"="17.399086"; "O"="0.2811226515149324"; "P"="0.2811226515149324"; "Q"="388.7157603132964"; "R"="3498.441842819668"; "S"="0.004113858231203587"; "T"="0.004113858231203587" }
    16 = @{ "E"="3"; "G"="67.02347933333333"; "H"="201.070438"; "I"="0.01463367753909034"; "J"="0.01463367753909034"; "K"="3"; "M"="6.259039333333333"; "N"="18.777118"; "O"="0.303387959572633"; "P"="0.303387959572633"; "Q"="419.5025934041871"; "R"="3775.523340637683"; "S"="0.004439681569628487"; "T"="0.004439681569628488" }
    17 = @{ "E"="3"; "G"="67.02347933333333"; "H"="201.070438"; "I"="0.01463367753909034"; "J"="0.01463367753909034"; "K"="3"; "M"="5.549851333333334"; "N"="16.649554"; "O"="0.2690122209305161"; "P"="0.2690122209305161"; "Q"="371.9703461427391"; "R"="3347.733115284652"; "S"="0.003936638095171701"; "T"="0.003936638095171701" }
}

$ws = $excel.ActiveWorkbook.ActiveSheet

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = [double]$cols[$col]
    }
}
